$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": updated aggregate metrics after trade #44 closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.75   # Current Capital
$summary.Range("B4").Value = 0.75      # Total P&L $
$summary.Range("B5").Value = 0.34      # Total P&L %
$summary.Range("B6").Value = 44        # Total Trades
$summary.Range("B8").Value = 21        # Losing Trades
$summary.Range("B9").Value = 31.82     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": MarketMaking row (row 4) refreshed to match.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.75     # Capital
$status.Range("D4").Value = 44         # Trades
$status.Range("E4").Value = 0.75       # P&L $
$status.Range("F4").Value = 0.75       # P&L %
$status.Range("G4").Value = 31.82      # Win Rate %

# ---------------------------------------------------------------------------
# New trade #44 (row 45) appended to both "All Trades" and "MarketMaking"
# sheets - identical data in each.
# ---------------------------------------------------------------------------
$newRow = @{
    A = 44
    B = "2026-02-17"
    C = "15:30:16"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.4
    G = 0.38
    H = "CLOSED"
    I = -5
    J = -0.02
    K = 100.75
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(45, 1).Value = $newRow.A

    # B45 looks like a date ("2026-02-17") - writing it with .Value causes
    # automatic date-serial conversion, so force it through as literal text
    # and then strip the format change back off so no new cell style is
    # introduced (matches every other date-like cell in these sheets, which
    # are stored as plain text).
    $ws.Cells.Item(45, 2).NumberFormat = "@"
    $ws.Cells.Item(45, 2).Value = $newRow.B
    $ws.Cells.Item(45, 2).ClearFormats()

    $ws.Cells.Item(45, 3).Value = $newRow.C
    $ws.Cells.Item(45, 4).Value = $newRow.D
    $ws.Cells.Item(45, 5).Value = $newRow.E
    $ws.Cells.Item(45, 6).Value = $newRow.F
    $ws.Cells.Item(45, 7).Value = $newRow.G
    $ws.Cells.Item(45, 8).Value = $newRow.H
    $ws.Cells.Item(45, 9).Value = $newRow.I
    $ws.Cells.Item(45, 10).Value = $newRow.J
    $ws.Cells.Item(45, 11).Value = $newRow.K
    $ws.Cells.Item(45, 12).Value = $newRow.L
    $ws.Cells.Item(45, 13).Value = $newRow.M
    $ws.Cells.Item(45, 14).Value = $newRow.N
    $ws.Cells.Item(45, 15).Value = $newRow.O
    $ws.Cells.Item(45, 16).Value = $newRow.P
    $ws.Cells.Item(45, 17).Value = $newRow.Q
}
